$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.054472091178195
$ws.Range("D2").Value = 1.052062178747892
$ws.Range("E2").Value = 1.060214665458465
$ws.Range("F2").Value = 1.069049703663304
$ws.Range("I2").Value = 1.044997968271764
$ws.Range("J2").Value = 1.059484010098365
$ws.Range("K2").Value = 1.05481195174809
$ws.Range("L2").Value = 1.062942053370013
$ws.Range("M2").Value = 1.071753249373307
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05565844617706
$ws.Range("D3").Value = 1.052942796986936
$ws.Range("E3").Value = 1.061308740310811
$ws.Range("F3").Value = 1.070313874540433
$ws.Range("I3").Value = 1.045327820263063
$ws.Range("J3").Value = 1.06032051276505
$ws.Range("K3").Value = 1.055505166449623
$ws.Range("L3").Value = 1.063849803412115
$ws.Range("M3").Value = 1.07283240419349
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056425902941884
$ws.Range("D4").Value = 1.053512372662861
$ws.Range("E4").Value = 1.062016841348914
$ws.Range("F4").Value = 1.071132351493372
$ws.Range("I4").Value = 1.045539949128471
$ws.Range("J4").Value = 1.060861044696932
$ws.Range("K4").Value = 1.055952842410646
$ws.Range("L4").Value = 1.064436744940241
$ws.Range("M4").Value = 1.073530600518337
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056748496868093
$ws.Range("D5").Value = 1.053751764486733
$ws.Range("E5").Value = 1.062314566614622
$ws.Range("F5").Value = 1.071476553976798
$ws.Range("I5").Value = 1.045628815433396
$ws.Range("J5").Value = 1.061088107726274
$ws.Range("K5").Value = 1.05614083521901
$ws.Range("L5").Value = 1.064683392279686
$ws.Range("M5").Value = 1.073824102074841
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056802659221685
$ws.Range("D6").Value = 1.053791956051097
$ws.Range("E6").Value = 1.0623645583456
$ws.Range("F6").Value = 1.071534353857305
$ws.Range("I6").Value = 1.045643718151541
$ws.Range("J6").Value = 1.061126222277336
$ws.Range("K6").Value = 1.056172387731416
$ws.Range("L6").Value = 1.064724799446275
$ws.Range("M6").Value = 1.073873381127232
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056430213632116
$ws.Range("D7").Value = 1.053515571654664
$ws.Range("E7").Value = 1.062020819411752
$ws.Range("F7").Value = 1.071136950291156
$ws.Range("I7").Value = 1.045541137792551
$ws.Range("J7").Value = 1.06086407941707
$ws.Range("K7").Value = 1.055955355204835
$ws.Range("L7").Value = 1.064440041057194
$ws.Range("M7").Value = 1.073534522378372
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05487306659136
$ws.Range("D8").Value = 1.052359838425652
$ws.Range("E8").Value = 1.060584380312813
$ws.Range("F8").Value = 1.069476838813296
$ws.Range("I8").Value = 1.045109714207894
$ws.Range("J8").Value = 1.059766864083431
$ws.Range("K8").Value = 1.055046409525278
$ws.Range("L8").Value = 1.063248922030868
$ws.Range("M8").Value = 1.072117973906297
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05212761472077
$ws.Range("D9").Value = 1.050321406101123
$ws.Range("E9").Value = 1.058054385209999
$ws.Range("F9").Value = 1.066555079713013
$ws.Range("I9").Value = 1.044339462350879
$ws.Range("J9").Value = 1.057827714837714
$ws.Range("K9").Value = 1.053437962622458
$ws.Range("L9").Value = 1.061146658248612
$ws.Range("M9").Value = 1.069621101950413
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050296161562438
$ws.Range("D10").Value = 1.048961154695234
$ws.Range("E10").Value = 1.056368462766406
$ws.Range("F10").Value = 1.06460955401516
$ws.Range("I10").Value = 1.043819196306214
$ws.Range("J10").Value = 1.056531049432968
$ws.Range("K10").Value = 1.052361074897676
$ws.Range("L10").Value = 1.059742831601926
$ws.Range("M10").Value = 1.067955954523955
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049502825210934
$ws.Range("D11").Value = 1.048371835619489
$ws.Range("E11").Value = 1.055638600882923
$ws.Range("F11").Value = 1.063767648700039
$ws.Range("I11").Value = 1.043592305466183
$ws.Range("J11").Value = 1.05596864089425
$ws.Range("K11").Value = 1.051893672619785
$ws.Range("L11").Value = 1.05913439575592
$ws.Range("M11").Value = 1.067234775882555
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049208096796691
$ws.Range("D12").Value = 1.048152887250181
$ws.Range("E12").Value = 1.055367519411918
$ws.Range("F12").Value = 1.063455003907305
$ws.Range("I12").Value = 1.043507785362173
$ws.Range("J12").Value = 1.055759594280528
$ws.Range("K12").Value = 1.051719892024483
$ws.Range("L12").Value = 1.058908308732193
$ws.Range("M12").Value = 1.066966872589969
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049271319247161
$ws.Range("D13").Value = 1.048199854649415
$ws.Range("E13").Value = 1.055425666346158
$ws.Range("F13").Value = 1.063522063851931
$ws.Range("I13").Value = 1.043525926204532
$ws.Range("J13").Value = 1.055804441981598
$ws.Range("K13").Value = 1.051757176104533
$ws.Range("L13").Value = 1.058956809117387
$ws.Range("M13").Value = 1.067024339931927
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04947846385923
$ws.Range("D14").Value = 1.048353738283503
$ws.Range("E14").Value = 1.055616192768542
$ws.Range("F14").Value = 1.063741803837424
$ws.Range("I14").Value = 1.043585323962968
$ws.Range("J14").Value = 1.055951363961105
$ws.Range("K14").Value = 1.051879311264907
$ws.Range("L14").Value = 1.059115709111323
$ws.Range("M14").Value = 1.067212631440789
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04960608607082
$ws.Range("D15").Value = 1.048448544550085
$ws.Range("E15").Value = 1.055733585224309
$ws.Range("F15").Value = 1.063877202954435
$ws.Range("I15").Value = 1.043621888700637
$ws.Range("J15").Value = 1.056041868450338
$ws.Range("K15").Value = 1.051954540668082
$ws.Range("L15").Value = 1.059213601093918
$ws.Range("M15").Value = 1.067328640636333
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05034880607368
$ws.Range("D16").Value = 1.049000259028265
$ws.Range("E16").Value = 1.056416904417636
$ws.Range("F16").Value = 1.064665439301898
$ws.Range("I16").Value = 1.043834220291604
$ws.Range("J16").Value = 1.056568354643853
$ws.Range("K16").Value = 1.052392071536554
$ws.Range("L16").Value = 1.059783199374473
$ws.Range("M16").Value = 1.068003813294831
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050814611529738
$ws.Range("D17").Value = 1.049346248402562
$ws.Range("E17").Value = 1.056845572677899
$ws.Range("F17").Value = 1.06516001711371
$ws.Range("I17").Value = 1.043966978282587
$ws.Range("J17").Value = 1.056898351953081
$ws.Range("K17").Value = 1.052666227211552
$ws.Range("L17").Value = 1.060140339916893
$ws.Range("M17").Value = 1.068427287635779
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05108627868243
$ws.Range("D18").Value = 1.049548027207232
$ws.Range("E18").Value = 1.057095622538742
$ws.Range("F18").Value = 1.065448546292965
$ws.Range("I18").Value = 1.044044258253206
$ws.Range("J18").Value = 1.057090742814565
$ws.Range("K18").Value = 1.05282603123582
$ws.Range("L18").Value = 1.060348599076379
$ws.Range("M18").Value = 1.0686742778364
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051178905258112
$ws.Range("D19").Value = 1.049616823339729
$ws.Range("E19").Value = 1.057180885650819
$ws.Range("F19").Value = 1.065546935875446
$ws.Range("I19").Value = 1.044070582363338
$ws.Range("J19").Value = 1.057156327731298
$ws.Range("K19").Value = 1.052880502268083
$ws.Range("L19").Value = 1.060419600798659
$ws.Range("M19").Value = 1.068758492578646
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050764638057318
$ws.Range("D20").Value = 1.049309130223091
$ws.Range("E20").Value = 1.056799579097071
$ws.Range("F20").Value = 1.065106948404302
$ws.Range("I20").Value = 1.043952750696382
$ws.Range("J20").Value = 1.056862955784693
$ws.Range("K20").Value = 1.052636823908036
$ws.Range("L20").Value = 1.060102027801796
$ws.Range("M20").Value = 1.068381854431901
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049417466275324
$ws.Range("D21").Value = 1.04830842474819
$ws.Range("E21").Value = 1.055560086896959
$ws.Range("F21").Value = 1.063677093796196
$ws.Range("I21").Value = 1.043567839508711
$ws.Range("J21").Value = 1.055908103058655
$ws.Range("K21").Value = 1.051843350109905
$ws.Range("L21").Value = 1.059068919431616
$ws.Range("M21").Value = 1.067157184999248
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048570167306032
$ws.Range("D22").Value = 1.047678957473418
$ws.Range("E22").Value = 1.054780894002622
$ws.Range("F22").Value = 1.062778528155717
$ws.Range("I22").Value = 1.043324425866371
$ws.Range("J22").Value = 1.055306921045502
$ws.Range("K22").Value = 1.051343497424134
$ws.Range("L22").Value = 1.058418860016443
$ws.Range("M22").Value = 1.066387038781508
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04901936364468
$ws.Range("D23").Value = 1.048012677163748
$ws.Range("E23").Value = 1.055193947499534
$ws.Range("F23").Value = 1.063254833638199
$ws.Range("I23").Value = 1.043453597369475
$ws.Range("J23").Value = 1.055625697923546
$ws.Range("K23").Value = 1.051608570487221
$ws.Range("L23").Value = 1.058763516934851
$ws.Range("M23").Value = 1.066795322412529
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050787219008352
$ws.Range("D24").Value = 1.049325902425991
$ws.Range("E24").Value = 1.056820361569979
$ws.Range("F24").Value = 1.065130927713814
$ws.Range("I24").Value = 1.043959180010522
$ws.Range("J24").Value = 1.056878950070337
$ws.Range("K24").Value = 1.052650110323533
$ws.Range("L24").Value = 1.06011933956666
$ws.Range("M24").Value = 1.068402383787652
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052837575302803
$ws.Range("D25").Value = 1.050848615973955
$ws.Range("E25").Value = 1.058708314441936
$ws.Range("F25").Value = 1.06731001034569
$ws.Range("I25").Value = 1.044539781627035
$ws.Range("J25").Value = 1.058329714527376
$ws.Range("K25").Value = 1.053854590812246
$ws.Range("L25").Value = 1.061690548001803
$ws.Range("M25").Value = 1.07026669758315
